$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the class name and start time on row 2
$ws.Range("A2").Value = "CICLO INDOOR"
$ws.Range("C2").Value = "19:00"

# Move the active selection to C7
$ws.Range("C7").Select()
